$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.035.31"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.681.79"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0620"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "1.919.91"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "1.687.87"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "27.057.62"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -3.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "1.512.96"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.589"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.919"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  +7.65%  "
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "1.823.85"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.70%  "
